$d = $word.ActiveDocument

# --- 1. Text fixes: swap which wrist goes with which direction -------------

$d.Content.Find.Execute(
    "right wrist above the left wrist", $true, $false, $false, $false,
    $false, $true, 1, $false, "left wrist above the right wrist", 2
) | Out-Null

$d.Content.Find.Execute(
    "Toggle Modes: right hand to the right, quickly", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Toggle Modes: left hand to the left, quickly", 2
) | Out-Null

$d.Content.Find.Execute(
    "Left: right wrist above left wrist", $true, $false, $false, $false,
    $false, $true, 1, $false, "Left: left wrist above right wrist", 2
) | Out-Null

$d.Content.Find.Execute(
    "Right: left wrist above right wrist", $true, $false, $false, $false,
    $false, $true, 1, $false, "Right: right wrist above left wrist", 2
) | Out-Null

# --- 2. Insert a new "Stop:" bullet right after the "Right: ..." bullet ---

$rightPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Right: right wrist above left wrist`r") {
        $rightPara = $d.Paragraphs.Item($i)
        break
    }
}

$rightPara.Range.InsertParagraphAfter()
$stopPara = $rightPara.Next()
$stopPara.Range.Text = "Stop: both wrists near both elbows in height"

# --- 3. Insert the new Arm Mode sub-bullets after "Arm Mode:" --------------

$armPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Arm Mode:`r") {
        $armPara = $d.Paragraphs.Item($i)
        break
    }
}

$newItems = @(
    "Hand Open: right hand to the right, medium quickness",
    "Hand Close: right hand to the left, medium quickness",
    "Wrist Up: Not defined yet",
    "Wrist Down: Not defined yet",
    "Shoulder Left: Not defined yet",
    "Shoulder Right: Not defined yet",
    "Shoulder Up: Not defined yet",
    "Shoulder Down: Not defined yet",
    "Elbow Up: Not defined yet",
    "Elbow Down: Not defined yet"
)

$cur = $armPara
foreach ($item in $newItems) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    $cur.Range.Text = $item
    # match the indentation level of the sibling bullets (Forward/Backward/...)
    $cur.Range.ListFormat.ListLevelNumber = 3
}

Write-Output "edit applied"
